# Scheduled-runner refresh of the profit/price figures on each job-sheet
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Updates the price/profit
# columns (H, I, J, K, L, M, N) for the rows whose market data changed;
# everything else on the sheets is left untouched.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H18").Value = 747.3158
$ws_ALC.Range("I18").Value = 622.1667
$ws_ALC.Range("K18").Value = 622.1667
$ws_ALC.Range("M18").Value = -338.1667

$ws_ALC.Range("H40").Value = 1648.1724
$ws_ALC.Range("I40").Value = 1540.7727
$ws_ALC.Range("K40").Value = 1540.7727
$ws_ALC.Range("M40").Value = -1365.7727

$ws_ALC.Range("H51").Value = 7184.1577
$ws_ALC.Range("J51").Value = 2051
$ws_ALC.Range("L51").Value = 2051
$ws_ALC.Range("N51").Value = -3019

$ws_ALC.Range("H76").Value = 3390.9092
$ws_ALC.Range("I76").Value = 3062.5
$ws_ALC.Range("J76").Value = 4266.6665
$ws_ALC.Range("K76").Value = 3062.5
$ws_ALC.Range("L76").Value = 4266.6665
$ws_ALC.Range("M76").Value = -2747.5
$ws_ALC.Range("N76").Value = -4896.6665

$ws_ALC.Range("H79").Value = 3390.9092
$ws_ALC.Range("I79").Value = 3062.5
$ws_ALC.Range("J79").Value = 4266.6665
$ws_ALC.Range("K79").Value = 3062.5
$ws_ALC.Range("L79").Value = 4266.6665
$ws_ALC.Range("M79").Value = -1970.5
$ws_ALC.Range("N79").Value = -6450.6665

$ws_ALC.Range("H92").Value = 714.26666
$ws_ALC.Range("I92").Value = 756.46155
$ws_ALC.Range("K92").Value = 756.46155
$ws_ALC.Range("M92").Value = 491.53845

$ws_ALC.Range("H96").Value = 667.7143
$ws_ALC.Range("I96").Value = 481.54544
$ws_ALC.Range("J96").Value = 1350.3334
$ws_ALC.Range("K96").Value = 1444.63632
$ws_ALC.Range("L96").Value = 4051.0002
$ws_ALC.Range("M96").Value = -71.63632000000007
$ws_ALC.Range("N96").Value = -6797.0002

$ws_ALC.Range("H132").Value = 4549804.5
$ws_ALC.Range("I132").Value = 5439594.5
$ws_ALC.Range("J132").Value = 1988.8889
$ws_ALC.Range("K132").Value = 16318783.5
$ws_ALC.Range("L132").Value = 5966.6667
$ws_ALC.Range("M132").Value = -16316253.5
$ws_ALC.Range("N132").Value = -11026.6667

$ws_ALC.Range("H137").Value = 1818.0923
$ws_ALC.Range("I137").Value = 1280.3928
$ws_ALC.Range("J137").Value = 2225
$ws_ALC.Range("K137").Value = 3841.1784
$ws_ALC.Range("L137").Value = 6675
$ws_ALC.Range("M137").Value = -1291.1784
$ws_ALC.Range("N137").Value = -11775

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 25225.404
$ws_ARM.Range("I2").Value = 1096.6451
$ws_ARM.Range("J2").Value = 93224.63
$ws_ARM.Range("K2").Value = 1096.6451
$ws_ARM.Range("L2").Value = 93224.63
$ws_ARM.Range("M2").Value = -983.6451
$ws_ARM.Range("N2").Value = -93450.63

$ws_ARM.Range("H32").Value = 675.13
$ws_ARM.Range("I32").Value = 614.7368
$ws_ARM.Range("J32").Value = 1822.6
$ws_ARM.Range("K32").Value = 614.7368
$ws_ARM.Range("L32").Value = 1822.6
$ws_ARM.Range("M32").Value = -327.7368
$ws_ARM.Range("N32").Value = -2396.6

$ws_ARM.Range("H116").Value = 25225.404
$ws_ARM.Range("I116").Value = 1096.6451
$ws_ARM.Range("J116").Value = 93224.63
$ws_ARM.Range("K116").Value = 1096.6451
$ws_ARM.Range("L116").Value = 93224.63
$ws_ARM.Range("M116").Value = 1197.3549
$ws_ARM.Range("N116").Value = -97812.63

$ws_ARM.Range("H132").Value = 2636.6667
$ws_ARM.Range("I132").Value = 2987.756
$ws_ARM.Range("J132").Value = 1879.0526
$ws_ARM.Range("K132").Value = 8963.268
$ws_ARM.Range("L132").Value = 5637.1578
$ws_ARM.Range("M132").Value = -6433.268
$ws_ARM.Range("N132").Value = -10697.1578

$ws_ARM.Range("H139").Value = 46131.11
$ws_ARM.Range("J139").Value = 46131.11
$ws_ARM.Range("L139").Value = 46131.11
$ws_ARM.Range("N139").Value = -56411.11

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 25225.404
$ws_BSM.Range("I3").Value = 1096.6451
$ws_BSM.Range("J3").Value = 93224.63
$ws_BSM.Range("K3").Value = 1096.6451
$ws_BSM.Range("L3").Value = 93224.63
$ws_BSM.Range("M3").Value = -982.6451
$ws_BSM.Range("N3").Value = -93452.63

$ws_BSM.Range("H75").Value = 9803.5
$ws_BSM.Range("I75").Value = 9738
$ws_BSM.Range("K75").Value = 9738
$ws_BSM.Range("M75").Value = -8802

$ws_BSM.Range("H78").Value = 9803.5
$ws_BSM.Range("I78").Value = 9738
$ws_BSM.Range("K78").Value = 29214
$ws_BSM.Range("M78").Value = -24534

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H7").Value = 131.17647
$ws_CRP.Range("I7").Value = 39.727272
$ws_CRP.Range("J7").Value = 298.83334
$ws_CRP.Range("K7").Value = 39.727272
$ws_CRP.Range("L7").Value = 298.83334
$ws_CRP.Range("M7").Value = 73.272728
$ws_CRP.Range("N7").Value = -524.83334

$ws_CRP.Range("H22").Value = 789
$ws_CRP.Range("I22").Value = 314.66666
$ws_CRP.Range("K22").Value = 314.66666
$ws_CRP.Range("M22").Value = 35.33334000000002

$ws_CRP.Range("H31").Value = 18648.262
$ws_CRP.Range("I31").Value = 44712.957
$ws_CRP.Range("J31").Value = 2872.2632
$ws_CRP.Range("K31").Value = 44712.957
$ws_CRP.Range("L31").Value = 2872.2632
$ws_CRP.Range("M31").Value = -44417.957
$ws_CRP.Range("N31").Value = -3462.2632

$ws_CRP.Range("H34").Value = 18648.262
$ws_CRP.Range("I34").Value = 44712.957
$ws_CRP.Range("J34").Value = 2872.2632
$ws_CRP.Range("K34").Value = 44712.957
$ws_CRP.Range("L34").Value = 2872.2632
$ws_CRP.Range("M34").Value = -44510.957
$ws_CRP.Range("N34").Value = -3276.2632

$ws_CRP.Range("H134").Value = 2175.0715
$ws_CRP.Range("I134").Value = 1044.7
$ws_CRP.Range("J134").Value = 5001
$ws_CRP.Range("K134").Value = 3134.1
$ws_CRP.Range("L134").Value = 15003
$ws_CRP.Range("M134").Value = -599.1000000000004
$ws_CRP.Range("N134").Value = -20073

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H19").Value = 1250
$ws_CUL.Range("J19").Value = 1250
$ws_CUL.Range("L19").Value = 3750
$ws_CUL.Range("N19").Value = -4098

$ws_CUL.Range("H92").Value = 1700
$ws_CUL.Range("I92").Value = 1500
$ws_CUL.Range("K92").Value = 4500
$ws_CUL.Range("M92").Value = -3252

$ws_CUL.Range("H113").Value = 558.0476
$ws_CUL.Range("I113").Value = 529.125
$ws_CUL.Range("J113").Value = 575.8461
$ws_CUL.Range("K113").Value = 1587.375
$ws_CUL.Range("L113").Value = 1727.5383
$ws_CUL.Range("M113").Value = 582.625
$ws_CUL.Range("N113").Value = -6067.5383

$ws_CUL.Range("H131").Value = 1390282.1
$ws_CUL.Range("J131").Value = 1588816.2
$ws_CUL.Range("L131").Value = 4766448.6
$ws_CUL.Range("N131").Value = -4776528.6

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 3755.8
$ws_GSM.Range("I80").Value = 4144.75
$ws_GSM.Range("J80").Value = 2200
$ws_GSM.Range("K80").Value = 4144.75
$ws_GSM.Range("L80").Value = 2200
$ws_GSM.Range("M80").Value = -3146.75
$ws_GSM.Range("N80").Value = -4196

$ws_GSM.Range("H83").Value = 3755.8
$ws_GSM.Range("I83").Value = 4144.75
$ws_GSM.Range("J83").Value = 2200
$ws_GSM.Range("K83").Value = 20723.75
$ws_GSM.Range("L83").Value = 11000
$ws_GSM.Range("M83").Value = -15731.75
$ws_GSM.Range("N83").Value = -20984

$ws_GSM.Range("H134").Value = 1454634.6
$ws_GSM.Range("J134").Value = 1454634.6
$ws_GSM.Range("L134").Value = 4363903.800000001
$ws_GSM.Range("N134").Value = -4368973.800000001

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 1809.6666
$ws_LTW.Range("I7").Value = 1371
$ws_LTW.Range("J7").Value = 2248.3333
$ws_LTW.Range("K7").Value = 1371
$ws_LTW.Range("L7").Value = 2248.3333
$ws_LTW.Range("M7").Value = -1259
$ws_LTW.Range("N7").Value = -2472.3333

$ws_LTW.Range("H46").Value = 1447287.1
$ws_LTW.Range("J46").Value = 2531752.5
$ws_LTW.Range("L46").Value = 2531752.5
$ws_LTW.Range("N46").Value = -2532128.5

$ws_LTW.Range("H61").Value = 2249.6667
$ws_LTW.Range("I61").Value = 2099.6
$ws_LTW.Range("J61").Value = 3000
$ws_LTW.Range("K61").Value = 2099.6
$ws_LTW.Range("L61").Value = 3000
$ws_LTW.Range("M61").Value = -1897.6
$ws_LTW.Range("N61").Value = -3404

$ws_LTW.Range("H82").Value = 2225.9092
$ws_LTW.Range("I82").Value = 2195.7144
$ws_LTW.Range("J82").Value = 2278.75
$ws_LTW.Range("K82").Value = 2195.7144
$ws_LTW.Range("L82").Value = 2278.75
$ws_LTW.Range("M82").Value = -1834.7144
$ws_LTW.Range("N82").Value = -3000.75

$ws_LTW.Range("H85").Value = 2225.9092
$ws_LTW.Range("I85").Value = 2195.7144
$ws_LTW.Range("J85").Value = 2278.75
$ws_LTW.Range("K85").Value = 2195.7144
$ws_LTW.Range("L85").Value = 2278.75
$ws_LTW.Range("M85").Value = -947.7143999999998
$ws_LTW.Range("N85").Value = -4774.75

$ws_LTW.Range("H113").Value = 2249.6667
$ws_LTW.Range("I113").Value = 2099.6
$ws_LTW.Range("J113").Value = 3000
$ws_LTW.Range("K113").Value = 2099.6
$ws_LTW.Range("L113").Value = 3000
$ws_LTW.Range("M113").Value = 70.40000000000009
$ws_LTW.Range("N113").Value = -7340

$ws_LTW.Range("H126").Value = 1809.6666
$ws_LTW.Range("I126").Value = 1371
$ws_LTW.Range("J126").Value = 2248.3333
$ws_LTW.Range("K126").Value = 4113
$ws_LTW.Range("L126").Value = 6744.999899999999
$ws_LTW.Range("M126").Value = -1643
$ws_LTW.Range("N126").Value = -11684.9999

$ws_LTW.Range("H136").Value = 1421
$ws_LTW.Range("I136").Value = 1495.2069
$ws_LTW.Range("J136").Value = 1225.3636
$ws_LTW.Range("K136").Value = 4485.620699999999
$ws_LTW.Range("L136").Value = 3676.0908
$ws_LTW.Range("M136").Value = -1935.620699999999
$ws_LTW.Range("N136").Value = -8776.0908

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H54").Value = 6931.5454
$ws_WVR.Range("J54").Value = 6917.7
$ws_WVR.Range("L54").Value = 6917.7
$ws_WVR.Range("N54").Value = -7957.7

$ws_WVR.Range("H140").Value = 60548.285
$ws_WVR.Range("J140").Value = 60548.285
$ws_WVR.Range("L140").Value = 60548.285
$ws_WVR.Range("N140").Value = -70908.285

$ws_WVR.Range("H141").Value = 60000
$ws_WVR.Range("J141").Value = 60000
$ws_WVR.Range("L141").Value = 60000
$ws_WVR.Range("N141").Value = -70360
